$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Update "Recorded By" (column G) values: reorder the comma-separated
# list of recorders for specific rows (daily attendance processing).
$updates = @{
    2 = "backup@backdoor.com, System, system"
    3 = "System, dnasr281@gmail.com"
    6 = "System, dnasr281@gmail.com"
    7 = "System, admin@admin.com"
    10 = "System, dnasr281@gmail.com"
    11 = "System, dnasr281@gmail.com"
    12 = "System, dnasr281@gmail.com"
    13 = "System, dnasr281@gmail.com"
    14 = "System, dnasr281@gmail.com"
    15 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    18 = "System, dnasr281@gmail.com"
    19 = "System, dnasr281@gmail.com"
    20 = "System, dnasr281@gmail.com"
    21 = "System, dnasr281@gmail.com"
    22 = "System, dnasr281@gmail.com"
    24 = "System, dnasr281@gmail.com"
    28 = "backup@backdoor.com, System, system"
    29 = "System, dnasr281@gmail.com"
    32 = "System, dnasr281@gmail.com"
    33 = "System, admin@admin.com"
    36 = "System, dnasr281@gmail.com"
    37 = "System, dnasr281@gmail.com"
    38 = "System, dnasr281@gmail.com"
    39 = "System, dnasr281@gmail.com"
    40 = "System, dnasr281@gmail.com"
    41 = "System, dnasr281@gmail.com"
    43 = "System, dnasr281@gmail.com"
    44 = "System, dnasr281@gmail.com"
    45 = "System, dnasr281@gmail.com"
    46 = "System, dnasr281@gmail.com"
    47 = "System, dnasr281@gmail.com"
    48 = "System, dnasr281@gmail.com"
    50 = "System, dnasr281@gmail.com"
    54 = "backup@backdoor.com, System, system"
    55 = "System, dnasr281@gmail.com"
    58 = "System, dnasr281@gmail.com"
    59 = "System, admin@admin.com"
    62 = "System, dnasr281@gmail.com"
    63 = "System, dnasr281@gmail.com"
    64 = "System, dnasr281@gmail.com"
    65 = "System, dnasr281@gmail.com"
    66 = "System, dnasr281@gmail.com"
    67 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    70 = "System, dnasr281@gmail.com"
    71 = "System, dnasr281@gmail.com"
    72 = "System, dnasr281@gmail.com"
    73 = "System, dnasr281@gmail.com"
    74 = "System, dnasr281@gmail.com"
    76 = "System, dnasr281@gmail.com"
    83 = "System, dnasr281@gmail.com"
    84 = "System, dnasr281@gmail.com"
    85 = "System, dnasr281@gmail.com"
    86 = "System, dnasr281@gmail.com"
    90 = "System, dnasr281@gmail.com"
    92 = "System, dnasr281@gmail.com"
    93 = "System, dnasr281@gmail.com"
    94 = "System, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    99 = "System, dnasr281@gmail.com"
    101 = "System, dnasr281@gmail.com"
    109 = "System, dnasr281@gmail.com"
    110 = "System, dnasr281@gmail.com"
    111 = "System, dnasr281@gmail.com"
    112 = "System, dnasr281@gmail.com"
    116 = "System, dnasr281@gmail.com"
    118 = "System, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    120 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    127 = "System, dnasr281@gmail.com"
    135 = "System, dnasr281@gmail.com"
    136 = "System, dnasr281@gmail.com"
    137 = "System, dnasr281@gmail.com"
    138 = "System, dnasr281@gmail.com"
    142 = "System, dnasr281@gmail.com"
    144 = "System, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    146 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
    153 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $updates[$row]
}

